$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G4").Value = "Adding text to excel"
[void]$ws.Range("G4").Select()
